# Adding the weeks that each task should be done.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell "Weeks Assigned" - match the style of the other column headers
# (G1/I1, style index 7) by copying the format from G1 before setting the value.
[void]$ws.Range("G1").Copy()
[void]$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("L1").Value = "Weeks Assigned"

# First week value is center-aligned (new style), the rest of the column
# reuses the default (no explicit style) the way the source workbook does.
$ws.Range("L2").Value = "Week 1-3"
$ws.Range("L2").HorizontalAlignment = -4108

$ws.Range("L3").Value = "Week 1-3"
$ws.Range("L4").Value = "Week 1-3"

$ws.Range("L10").Value = "Week 2-3"
$ws.Range("L11").Value = "Week 2-3"

$ws.Range("L9").Value = "Week 1-2"
# Stray centered-but-empty formatting artifact in O9
$ws.Range("O9").HorizontalAlignment = -4108

$ws.Range("L12").Value = "Week 2-4"
$ws.Range("L15").Value = "Week 2-7"
$ws.Range("L26").Value = "Week 3-9"
$ws.Range("L37").Value = "Week 5-11"
$ws.Range("L63").Value = "Week 12"

# Restore the selection to where editing finished
[void]$ws.Range("L63").Select()
